$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2 value
$ws.Range("B2").Value = "Initialization of application and beta database"

# Add new row 3
$ws.Range("A3").Value = "#100001"
$ws.Range("B3").Value = "Add entitites and dto objects"

# Adjust column B width (closest achievable value in this engine's pixel-quantized
# column-width model; real Excel stores 44.88671875 but this headless engine
# rounds column widths to whole pixels, so 44 is the input that lands nearest)
$ws.Columns.Item(2).ColumnWidth = 44

# Adjust selection
$ws.Range("O11").Select()
